# Events.xlsx update: add "Omikron prevalent variant" event row and RKI source entry.

$wb = $excel.ActiveWorkbook

# --- Sheet1: Events table -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at position 39 (shifts existing rows 39-41 down to 40-42)
$ws1.Rows.Item(39).Insert()

# Fill in the new event row
$ws1.Range("A39").Value = 44557
$ws1.Range("B39").Value = "Omikron prevalent variant"
$ws1.Range("C39").Value = "RKI"
$ws1.Range("D39").Value = "x"

# Update selection to match the saved workbook state
$ws1.Activate()
$ws1.Range("A40").Select()

# --- Sheet2: Sources / legend table ---------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Expand the "MDR" abbreviation entries to the full name
$ws2.Range("B2").Value = "Mitteldeutscher Rundfunk"
$ws2.Range("B3").Value = "Mitteldeutscher Rundfunk"
$ws2.Range("B4").Value = "Mitteldeutscher Rundfunk"

# Add a new row describing the RKI source
$ws2.Range("A6").Value = "RKI"
$ws2.Range("B6").Value = "Robert-Koch-Institut"
$ws2.Range("C6").Value = "https://www.rki.de/DE/Content/InfAZ/N/Neuartiges_Coronavirus/Virusvariante.html"
$ws2.Range("D6").Value = "Table under the text"

# Update selection to match the saved workbook state
$ws2.Activate()
$ws2.Range("B7").Select()
